$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet (positioned right after 2021-Q4,
#    right before the existing total/summary sheet). We copy an
#    existing quarter sheet so we inherit its sheetPr/format boiler-
#    plate (outline props, page margins, etc.), then wipe its cells.
# ---------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($null, $templateSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet.Cells.Clear()

# Re-apply the header-row / index-column styles by copying them from
# the template sheet (keeps the same style indices: bold+border for
# the header row, bold+border for column A).
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row text
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows (index col A, fund code B, fund name C, fund size D,
# stock allocation E, position ratio F, holding value G, rank H)
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'010659"
$newSheet.Cells.Item(2,4).Value = "'24.14"
$newSheet.Cells.Item(2,5).Value = "'87.59"
$newSheet.Cells.Item(2,6).Value = "'3.98"
$newSheet.Cells.Item(2,7).Value = "'0.9608"
$newSheet.Cells.Item(2,3).Value = "民生加银质量领先混合A"
$newSheet.Cells.Item(2,8).Value = 10

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'009362"
$newSheet.Cells.Item(3,4).Value = "'21.39"
$newSheet.Cells.Item(3,5).Value = "'87.04"
$newSheet.Cells.Item(3,6).Value = "'4.35"
$newSheet.Cells.Item(3,7).Value = "'0.9305"
$newSheet.Cells.Item(3,3).Value = "招商丰盈积极配置混合A"
$newSheet.Cells.Item(3,8).Value = 6

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'952004"
$newSheet.Cells.Item(4,4).Value = "'22.12"
$newSheet.Cells.Item(4,5).Value = "'76.15"
$newSheet.Cells.Item(4,6).Value = "'2.54"
$newSheet.Cells.Item(4,7).Value = "'0.5618"
$newSheet.Cells.Item(4,3).Value = "国泰君安君得明混合"
$newSheet.Cells.Item(4,8).Value = 7

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'217010"
$newSheet.Cells.Item(5,4).Value = "'9.19"
$newSheet.Cells.Item(5,5).Value = "'84.38"
$newSheet.Cells.Item(5,6).Value = "'5.35"
$newSheet.Cells.Item(5,7).Value = "'0.4917"
$newSheet.Cells.Item(5,3).Value = "招商大盘蓝筹混合"
$newSheet.Cells.Item(5,8).Value = 4

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'217009"
$newSheet.Cells.Item(6,4).Value = "'8.46"
$newSheet.Cells.Item(6,5).Value = "'82.51"
$newSheet.Cells.Item(6,6).Value = "'4.90"
$newSheet.Cells.Item(6,7).Value = "'0.4145"
$newSheet.Cells.Item(6,3).Value = "招商核心价值混合"
$newSheet.Cells.Item(6,8).Value = 4

$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'010795"
$newSheet.Cells.Item(7,4).Value = "'7.16"
$newSheet.Cells.Item(7,5).Value = "'88.23"
$newSheet.Cells.Item(7,6).Value = "'3.98"
$newSheet.Cells.Item(7,7).Value = "'0.2850"
$newSheet.Cells.Item(7,3).Value = "民生加银价值发现一年持有期混合A"
$newSheet.Cells.Item(7,8).Value = 10

$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'000408"
$newSheet.Cells.Item(8,4).Value = "'7.28"
$newSheet.Cells.Item(8,5).Value = "'78.80"
$newSheet.Cells.Item(8,6).Value = "'3.47"
$newSheet.Cells.Item(8,7).Value = "'0.2526"
$newSheet.Cells.Item(8,3).Value = "民生加银城镇化混合A"
$newSheet.Cells.Item(8,8).Value = 10

$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'011690"
$newSheet.Cells.Item(9,4).Value = "'7.06"
$newSheet.Cells.Item(9,5).Value = "'87.78"
$newSheet.Cells.Item(9,6).Value = "'3.47"
$newSheet.Cells.Item(9,7).Value = "'0.2450"
$newSheet.Cells.Item(9,3).Value = "招商品质发现混合A"
$newSheet.Cells.Item(9,8).Value = 5

$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'009695"
$newSheet.Cells.Item(10,4).Value = "'5.11"
$newSheet.Cells.Item(10,5).Value = "'87.04"
$newSheet.Cells.Item(10,6).Value = "'4.58"
$newSheet.Cells.Item(10,7).Value = "'0.2340"
$newSheet.Cells.Item(10,3).Value = "招商成长精选一年定期开放混合A"
$newSheet.Cells.Item(10,8).Value = 6

$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "'011843"
$newSheet.Cells.Item(11,4).Value = "'5.53"
$newSheet.Cells.Item(11,5).Value = "'80.72"
$newSheet.Cells.Item(11,6).Value = "'3.46"
$newSheet.Cells.Item(11,7).Value = "'0.1913"
$newSheet.Cells.Item(11,3).Value = "民生加银内核驱动混合型证券投资基金A"
$newSheet.Cells.Item(11,8).Value = 10

$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "'009363"
$newSheet.Cells.Item(12,4).Value = "'4.09"
$newSheet.Cells.Item(12,5).Value = "'87.04"
$newSheet.Cells.Item(12,6).Value = "'4.35"
$newSheet.Cells.Item(12,7).Value = "'0.1779"
$newSheet.Cells.Item(12,3).Value = "招商丰盈积极配置混合C"
$newSheet.Cells.Item(12,8).Value = 6

$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "'519156"
$newSheet.Cells.Item(13,4).Value = "'4.98"
$newSheet.Cells.Item(13,5).Value = "'93.77"
$newSheet.Cells.Item(13,6).Value = "'3.16"
$newSheet.Cells.Item(13,7).Value = "'0.1574"
$newSheet.Cells.Item(13,3).Value = "新华行业轮换灵活配置混合A"
$newSheet.Cells.Item(13,8).Value = 9

$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "'002249"
$newSheet.Cells.Item(14,4).Value = "'2.03"
$newSheet.Cells.Item(14,5).Value = "'87.69"
$newSheet.Cells.Item(14,6).Value = "'6.04"
$newSheet.Cells.Item(14,7).Value = "'0.1226"
$newSheet.Cells.Item(14,3).Value = "招商境远灵活配置混合"
$newSheet.Cells.Item(14,8).Value = 4

$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = "'011727"
$newSheet.Cells.Item(15,4).Value = "'3.46"
$newSheet.Cells.Item(15,5).Value = "'29.56"
$newSheet.Cells.Item(15,6).Value = "'2.06"
$newSheet.Cells.Item(15,7).Value = "'0.0713"
$newSheet.Cells.Item(15,3).Value = "工银瑞信聚瑞混合型证券投资基金A"
$newSheet.Cells.Item(15,8).Value = 3

$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = "'009696"
$newSheet.Cells.Item(16,4).Value = "'1.36"
$newSheet.Cells.Item(16,5).Value = "'87.04"
$newSheet.Cells.Item(16,6).Value = "'4.58"
$newSheet.Cells.Item(16,7).Value = "'0.0623"
$newSheet.Cells.Item(16,3).Value = "招商成长精选一年定期开放混合C"
$newSheet.Cells.Item(16,8).Value = 6

$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).Value = "'001707"
$newSheet.Cells.Item(17,4).Value = "'1.35"
$newSheet.Cells.Item(17,5).Value = "'91.74"
$newSheet.Cells.Item(17,6).Value = "'3.31"
$newSheet.Cells.Item(17,7).Value = "'0.0447"
$newSheet.Cells.Item(17,3).Value = "诺安高端制造股票"
$newSheet.Cells.Item(17,8).Value = 10

$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).Value = "'010660"
$newSheet.Cells.Item(18,4).Value = "'1.09"
$newSheet.Cells.Item(18,5).Value = "'87.59"
$newSheet.Cells.Item(18,6).Value = "'3.98"
$newSheet.Cells.Item(18,7).Value = "'0.0434"
$newSheet.Cells.Item(18,3).Value = "民生加银质量领先混合C"
$newSheet.Cells.Item(18,8).Value = 10

$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).Value = "'930602"
$newSheet.Cells.Item(19,4).Value = "'0.50"
$newSheet.Cells.Item(19,5).Value = "'67.38"
$newSheet.Cells.Item(19,6).Value = "'6.44"
$newSheet.Cells.Item(19,7).Value = "'0.0322"
$newSheet.Cells.Item(19,3).Value = "国信价值智选混合型集合资产管理计划"
$newSheet.Cells.Item(19,8).Value = 4

$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).Value = "'008518"
$newSheet.Cells.Item(20,4).Value = "'0.75"
$newSheet.Cells.Item(20,5).Value = "'81.65"
$newSheet.Cells.Item(20,6).Value = "'3.58"
$newSheet.Cells.Item(20,7).Value = "'0.0268"
$newSheet.Cells.Item(20,3).Value = "弘毅远方经济新动力混合"
$newSheet.Cells.Item(20,8).Value = 9

$newSheet.Cells.Item(21,1).Value = 19
$newSheet.Cells.Item(21,2).Value = "'010796"
$newSheet.Cells.Item(21,4).Value = "'0.52"
$newSheet.Cells.Item(21,5).Value = "'88.23"
$newSheet.Cells.Item(21,6).Value = "'3.98"
$newSheet.Cells.Item(21,7).Value = "'0.0207"
$newSheet.Cells.Item(21,3).Value = "民生加银价值发现一年持有期混合C"
$newSheet.Cells.Item(21,8).Value = 10

$newSheet.Cells.Item(22,1).Value = 20
$newSheet.Cells.Item(22,2).Value = "'010668"
$newSheet.Cells.Item(22,4).Value = "'1.03"
$newSheet.Cells.Item(22,5).Value = "'51.05"
$newSheet.Cells.Item(22,6).Value = "'1.11"
$newSheet.Cells.Item(22,7).Value = "'0.0114"
$newSheet.Cells.Item(22,3).Value = "工银瑞信优选对冲策略灵活配置混合A"
$newSheet.Cells.Item(22,8).Value = 4

$newSheet.Cells.Item(23,1).Value = 21
$newSheet.Cells.Item(23,2).Value = "'009706"
$newSheet.Cells.Item(23,4).Value = "'0.30"
$newSheet.Cells.Item(23,5).Value = "'78.80"
$newSheet.Cells.Item(23,6).Value = "'3.47"
$newSheet.Cells.Item(23,7).Value = "'0.0104"
$newSheet.Cells.Item(23,3).Value = "民生加银城镇化混合C"
$newSheet.Cells.Item(23,8).Value = 10

$newSheet.Cells.Item(24,1).Value = 22
$newSheet.Cells.Item(24,2).Value = "'010669"
$newSheet.Cells.Item(24,4).Value = "'0.62"
$newSheet.Cells.Item(24,5).Value = "'51.05"
$newSheet.Cells.Item(24,6).Value = "'1.11"
$newSheet.Cells.Item(24,7).Value = "'0.0069"
$newSheet.Cells.Item(24,3).Value = "工银瑞信优选对冲策略灵活配置混合C"
$newSheet.Cells.Item(24,8).Value = 4

$newSheet.Cells.Item(25,1).Value = 23
$newSheet.Cells.Item(25,2).Value = "'011844"
$newSheet.Cells.Item(25,4).Value = "'0.18"
$newSheet.Cells.Item(25,5).Value = "'80.72"
$newSheet.Cells.Item(25,6).Value = "'3.46"
$newSheet.Cells.Item(25,7).Value = "'0.0062"
$newSheet.Cells.Item(25,3).Value = "民生加银内核驱动混合型证券投资基金C"
$newSheet.Cells.Item(25,8).Value = 10

$newSheet.Cells.Item(26,1).Value = 24
$newSheet.Cells.Item(26,2).Value = "'011691"
$newSheet.Cells.Item(26,4).Value = "'0.13"
$newSheet.Cells.Item(26,5).Value = "'87.78"
$newSheet.Cells.Item(26,6).Value = "'3.47"
$newSheet.Cells.Item(26,7).Value = "'0.0045"
$newSheet.Cells.Item(26,3).Value = "招商品质发现混合C"
$newSheet.Cells.Item(26,8).Value = 5

$newSheet.Cells.Item(27,1).Value = 25
$newSheet.Cells.Item(27,2).Value = "'011728"
$newSheet.Cells.Item(27,4).Value = "'0.17"
$newSheet.Cells.Item(27,5).Value = "'29.56"
$newSheet.Cells.Item(27,6).Value = "'2.06"
$newSheet.Cells.Item(27,7).Value = "'0.0035"
$newSheet.Cells.Item(27,3).Value = "工银瑞信聚瑞混合型证券投资基金C"
$newSheet.Cells.Item(27,8).Value = 3

$newSheet.Cells.Item(28,1).Value = 26
$newSheet.Cells.Item(28,2).Value = "'519157"
$newSheet.Cells.Item(28,4).Value = "'0.04"
$newSheet.Cells.Item(28,5).Value = "'93.77"
$newSheet.Cells.Item(28,6).Value = "'3.16"
$newSheet.Cells.Item(28,7).Value = "'0.0013"
$newSheet.Cells.Item(28,3).Value = "新华行业轮换灵活配置混合C"
$newSheet.Cells.Item(28,8).Value = 9

Write-Host "2022-Q1 sheet populated"

# ---------------------------------------------------------------
# 2) Update the "总计" (total) summary sheet: insert a new row at
#    the top of the data for the 2022-Q1 quarter, pushing the
#    existing quarters down by one row.
# ---------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 27
$totalSheet.Cells.Item(2,4).Value = 5.37
Write-Host "Total sheet updated"

Write-Host "Sheet order:"
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
